# Resume edit: rewrite contact line, Education, Actuarial Exams, Experience,
# Projects and Skills sections.
#
# Strategy: the section below walks the document from the BOTTOM up,
# mutating each paragraph's Range.Text in place where only the wording
# changes, and inserting/deleting whole paragraphs where the bullet/line
# count of a block changes. Working bottom-to-top keeps every paragraph
# index we reference stable, since an insert/delete only renumbers the
# paragraphs that come AFTER it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# SKILLS -> TECHNICAL SKILLS (same number of bullets, wording only)
# ---------------------------------------------------------------------
$d.Paragraphs(34).Range.Text = "• Cloud Computing: Azure, AWS Sagemaker, Google Vertex AI"
$d.Paragraphs(33).Range.Text = "• Data Science & Machine Learning: Data Mining, NLP, Keras, TensorFlow, Pytorch, LangChain, Hugging Face"
$d.Paragraphs(32).Range.Text = "• Programming & Tools: Python, R, SQL, SAS, Microsoft Excel, Power BI, Tableau, Flask, Web App Development"
$d.Paragraphs(31).Range.Text = "TECHNICAL SKILLS"

# ---------------------------------------------------------------------
# Regression Modeling project (same number of bullets, wording only)
# ---------------------------------------------------------------------
$d.Paragraphs(29).Range.Text = "• Performed Exploratory Data Analysis (EDA), feature selection for accurate and effective predictions"
$d.Paragraphs(28).Range.Text = "• Optimized a multiple linear regression model using R for medical expenses prediction"
$d.Paragraphs(27).Range.Text = "Regression Modeling Project, California State University, East Bay, Spring 2024"

# ---------------------------------------------------------------------
# Traveler's case competition project (same number of bullets, wording only)
# ---------------------------------------------------------------------
$d.Paragraphs(25).Range.Text = "• Proposed reinsurance quotes and treaty structures enhancing the target market while ensuring risk tolerance alignment"
$d.Paragraphs(24).Range.Text = "• Utilized actuarial methods to assess future premiums and losses based on historical data"
$d.Paragraphs(23).Range.Text = "Reinsurance Analysis Project, Traveler’s 2025 Actuarial Case Competition, Spring 2025"

# ---------------------------------------------------------------------
# CAS summer program: 4 lines (title, dates, 2 bullets) -> 3 lines
# (title+dates merged, 2 reworded bullets) - drop the trailing bullet.
# ---------------------------------------------------------------------
$d.Paragraphs(20).Range.Delete()
$d.Paragraphs(19).Range.Text = "• Gained exposure to different areas such as data visualization, ratemaking, reserving, predictive modeling"
$d.Paragraphs(18).Range.Text = "• Participated in an eight-week interactive learning program focusing on property and casualty insurance"
$d.Paragraphs(17).Range.Text = "Summer Intern, Casualty Actuarial Society, June 2025 - August 2025"

# ---------------------------------------------------------------------
# Freelance experience: 4 lines (title, dates, 2 bullets) -> 3 lines
# (title+dates merged, 2 reworded bullets) - drop the trailing bullet.
# ---------------------------------------------------------------------
$d.Paragraphs(15).Range.Delete()
$d.Paragraphs(14).Range.Text = "• Provided insights on AI powered strategies for optimization"
$d.Paragraphs(13).Range.Text = "• Delivered statistical analysis, predictive modeling solutions for clients in various industries"
$d.Paragraphs(12).Range.Text = "Freelance Data Scientist, Fiverr.com/Upwork.com, July 2024 - Present"

# ---------------------------------------------------------------------
# CERTIFICATIONS -> ACTUARIAL EXAMS PASSED, turned into a bulleted list
# (old header + single summary line -> new header + 3 bullets: same
# total paragraph count, so pure in-place text rewrites).
# ---------------------------------------------------------------------
$d.Paragraphs(9).Range.Text = "• Fundamentals of Actuarial Mathematics (FAM)"
$d.Paragraphs(8).Range.Text = "• Financial Mathematics (FM)"
$d.Paragraphs(7).Range.Text = "• Probability (P)"
$d.Paragraphs(6).Range.Text = "ACTUARIAL EXAMS PASSED"

# ---------------------------------------------------------------------
# EDUCATION: the M.S./GPA line is split into its own line plus a new
# "GPA: 3.6" line and a new Bachelor's line (the old Bachelor's +
# trailing blank paragraph are consumed as part of this expansion).
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs(4)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "GPA: 3.6"
$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "Bachelor of Science in Business Administration, University of California, Berkeley"

$d.Paragraphs(4).Range.Text = "Master of Science in Statistics, Concentration in Actuarial Science, California State University, East Bay, Expected completion: May 2026 (part-time online)"

# ---------------------------------------------------------------------
# Header contact line: drop the trailing " |" and add a new blank
# paragraph after it (there wasn't one here before).
# ---------------------------------------------------------------------
$d.Paragraphs(2).Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "454 S Catalina St., Los Angeles, CA 90020 | (213) 433-4445"
